$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two Job Locations values: "charlotte, nc" (C2) and "ga" (C8)
$ws.Range("C2").Value = "ga"
$ws.Range("C8").Value = "charlotte, nc"

# Fix typo in Exact Phrases column: "on the job trainiing" -> "on the job training"
$ws.Range("B7").Value = "on the job training"

# Move the active cell selection to C2
$ws.Range("C2").Select()
